$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.853.26"
$ws.Range("E2").Value = "  -2.36%  "
$ws.Range("D3").Value = "3.387.91"
$ws.Range("E3").Value = "  -3.41%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "'574.40"
$ws.Range("E5").Value = "  -2.85%  "
$ws.Range("D6").Value = "'126.27"
$ws.Range("E6").Value = "  -6.16%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").Value = "3.387.27"
$ws.Range("E8").Value = "  -3.39%  "
$ws.Range("D9").Value = "'0.477"
$ws.Range("E9").Value = "  -2.19%  "
$ws.Range("D10").Value = "'7.27"
$ws.Range("E10").Value = "  -4.60%  "
$ws.Range("E11").Value = "  -4.64%  "
$ws.Range("D12").Value = "'0.377"
$ws.Range("E12").Value = "  -3.20%  "
$ws.Range("D13").Value = "3.964.31"
$ws.Range("E13").Value = "  -3.39%  "
$ws.Range("E14").Value = "  -0.97%  "
$ws.Range("D15").Value = "3.383.66"
$ws.Range("E15").Value = "  -3.18%  "
$ws.Range("E16").Value = "  -5.52%  "
$ws.Range("D17").Value = "62.903.63"
$ws.Range("E17").Value = "  -2.27%  "
$ws.Range("D18").Value = "'24.62"
$ws.Range("E18").Value = "  -4.59%  "
$ws.Range("D19").Value = "'9.31"
$ws.Range("E19").Value = "  -7.57%  "
$ws.Range("E20").Value = "  -2.32%  "
$ws.Range("D21").Value = "'13.13"
$ws.Range("E21").Value = "  -3.70%  "
$ws.Range("D22").Value = "'372.34"
$ws.Range("E22").Value = "  -5.31%  "
$ws.Range("D23").Value = "'0.555"
$ws.Range("E23").Value = "  -4.66%  "
$ws.Range("D24").Value = "3.520.54"
$ws.Range("E24").Value = "  -3.46%  "
$ws.Range("D25").Value = "'0.999"
$ws.Range("E25").Value = "  -0.38%  "
$ws.Range("D26").Value = "'71.57"
$ws.Range("E26").Value = "  -3.89%  "
$ws.Range("E27").Value = "  -9.95%  "
$ws.Range("D28").Value = "'0.999"
$ws.Range("E28").Value = "  -1.90%  "
$ws.Range("E29").Value = "  -5.84%  "
$ws.Range("B30").Value = "PancakeSwap"
$ws.Range("C30").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D30").Value = "'2.12"
$ws.Range("E30").Value = "  -6.81%  "
$ws.Range("B31").Value = "InternetComputer(DFINITY)"
$ws.Range("C31").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D31").Value = "'7.84"
$ws.Range("E31").Value = "  -4.96%  "
$ws.Range("E32").Value = "  -4.50%  "
$ws.Range("E33").Value = "  -0.02%  "
$ws.Range("E34").Value = "  -5.07%  "
$ws.Range("D35").Value = "3.416.37"
$ws.Range("E35").Value = "  -3.37%  "
$ws.Range("D36").Value = "'22.76"
$ws.Range("E36").Value = "  -2.75%  "
$ws.Range("D37").Value = "'5.36"
$ws.Range("E37").Value = "  +0.18%  "
$ws.Range("D38").Value = "'166.25"
$ws.Range("E38").Value = "  -0.01%  "
$ws.Range("D39").Value = "'6.66"
$ws.Range("E39").Value = "  -4.41%  "
$ws.Range("E40").Value = "  -4.80%  "
$ws.Range("D41").Value = "'0.0755"
$ws.Range("E41").Value = "  -4.23%  "
$ws.Range("D42").Value = "'1.00"
$ws.Range("E42").Value = "  +0.09%  "
$ws.Range("D43").Value = "'41.81"
$ws.Range("E43").Value = "  -0.56%  "
$ws.Range("D44").Value = "'0.763"
$ws.Range("E44").Value = "  -5.93%  "
$ws.Range("E45").Value = "  -4.83%  "
$ws.Range("D46").Value = "'1.55"
$ws.Range("E46").Value = "  -6.86%  "
$ws.Range("E47").Value = "  -6.85%  "
$ws.Range("D48").Value = "'22.38"
$ws.Range("E48").Value = "  -10.67%  "
$ws.Range("D49").Value = "'6.61"
$ws.Range("E49").Value = "  -2.89%  "
$ws.Range("D50").Value = "2.240.65"
$ws.Range("E50").Value = "  -5.88%  "
$ws.Range("E51").Value = "  -7.87%  "
